$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '29.772.88'
Set-TextValue 'E2' '  +2.09%  '
Set-TextValue 'D3' '1.859.86'
Set-TextValue 'E3' '  +1.67%  '
Set-TextValue 'D4' '0.9997'
Set-TextValue 'E4' '  +0.09%  '
Set-TextValue 'D5' '245.05'
Set-TextValue 'E5' '  +0.97%  '
Set-TextValue 'D6' '0.6407'
Set-TextValue 'E6' '  +3.32%  '
Set-TextValue 'B8' 'Dogecoin'
Set-TextValue 'C8' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D8' '0.07545'
Set-TextValue 'E8' '  +2.70%  '
Set-TextValue 'B9' 'Cardano'
Set-TextValue 'C9' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D9' '0.2977'
Set-TextValue 'E9' '  +2.67%  '
Set-TextValue 'B10' 'Solana'
Set-TextValue 'C10' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D10' '24.55'
Set-TextValue 'E10' '  +5.53%  '
Set-TextValue 'B11' 'TRON'
Set-TextValue 'C11' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D11' '0.07672'
Set-TextValue 'E11' '  +0.52%  '
Set-TextValue 'B12' 'WrappedEther'
Set-TextValue 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.859.50'
Set-TextValue 'E12' '  +1.56%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '5.042'
Set-TextValue 'E13' '  +1.61%  '
Set-TextValue 'B14' 'Polygon'
Set-TextValue 'C14' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.6920'
Set-TextValue 'E14' '  +3.47%  '
Set-TextValue 'B15' 'Litecoin'
Set-TextValue 'C15' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D15' '83.94'
Set-TextValue 'E15' '  +1.89%  '
Set-TextValue 'B16' 'ShibaInu'
Set-TextValue 'C16' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D16' '0.000009872'
Set-TextValue 'E16' '  +9.99%  '
Set-TextValue 'B17' 'Uniswap'
Set-TextValue 'C17' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D17' '6.112'
Set-TextValue 'E17' '  +4.75%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '29.786.19'
Set-TextValue 'E18' '  +2.18%  '
Set-TextValue 'B19' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C19' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D19' '2.113.49'
Set-TextValue 'E19' '  +1.57%  '
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '235.90'
Set-TextValue 'E20' '  +0.01%  '
Set-TextValue 'B21' 'Avalanche'
Set-TextValue 'C21' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D21' '12.66'
Set-TextValue 'E21' '  +1.66%  '
Set-TextValue 'B22' 'Dai'
Set-TextValue 'C22' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D22' '1.000'
Set-TextValue 'E22' '  +0.01%  '
Set-TextValue 'B23' 'Chainlink'
Set-TextValue 'C23' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D23' '7.505'
Set-TextValue 'E23' '  +2.17%  '
Set-TextValue 'B24' 'BinanceUSD'
Set-TextValue 'C24' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D24' '1.001'
Set-TextValue 'E24' '  +0.06%  '
Set-TextValue 'B25' 'Monero'
Set-TextValue 'C25' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D25' '159.25'
Set-TextValue 'E25' '  +0.58%  '
Set-TextValue 'B26' 'Stellar'
Set-TextValue 'C26' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D26' '0.1420'
Set-TextValue 'E26' '  +2.16%  '
Set-TextValue 'B27' 'Cosmos'
Set-TextValue 'C27' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D27' '8.550'
Set-TextValue 'E27' '  +0.39%  '
Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '17.94'
Set-TextValue 'E28' '  +1.85%  '
Set-TextValue 'B29' 'Hedera'
Set-TextValue 'C29' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D29' '0.06197'
Set-TextValue 'E29' '  +6.00%  '
Set-TextValue 'B30' 'PancakeSwap'
Set-TextValue 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '1.493'
Set-TextValue 'E30' '  +0.55%  '
Set-TextValue 'B31' 'Toncoin'
Set-TextValue 'C31' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '1.288'
Set-TextValue 'E31' '  +6.29%  '
Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '4.162'
Set-TextValue 'E32' '  +1.93%  '
Set-TextValue 'B33' 'InternetComputer(DFINITY)'
Set-TextValue 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D33' '4.104'
Set-TextValue 'E33' '  +0.68%  '
Set-TextValue 'B34' 'LidoDAOToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D34' '1.898'
Set-TextValue 'E34' '  +1.99%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.176'
Set-TextValue 'E35' '  +3.39%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.7290'
Set-TextValue 'E36' '  +0.59%  '
Set-TextValue 'B37' 'HuobiToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D37' '2.604'
Set-TextValue 'E37' '  -0.14%  '
Set-TextValue 'B38' 'MXToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D38' '2.827'
Set-TextValue 'E38' '  -0.62%  '
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.01786'
Set-TextValue 'E39' '  +1.63%  '
Set-TextValue 'B40' 'Maker'
Set-TextValue 'C40' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D40' '1.207.64'
Set-TextValue 'E40' '  -1.56%  '
Set-TextValue 'B41' 'FraxShare'
Set-TextValue 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D41' '6.288'
Set-TextValue 'E41' '  +1.09%  '
Set-TextValue 'D42' '0.9201'
Set-TextValue 'E42' '  +1.36%  '
Set-TextValue 'B43' 'PaxDollar'
Set-TextValue 'C43' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D43' '1.001'
Set-TextValue 'E43' '  +0.06%  '
Set-TextValue 'B44' 'RocketPoolETH'
Set-TextValue 'C44' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D44' '2.024.65'
Set-TextValue 'E44' '  +1.98%  '
Set-TextValue 'B45' 'Quant'
Set-TextValue 'C45' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D45' '102.04'
Set-TextValue 'E45' '  +0.28%  '
Set-TextValue 'B46' 'Aave'
Set-TextValue 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '66.56'
Set-TextValue 'E46' '  +1.29%  '
Set-TextValue 'B47' 'BabyDogeCoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D47' '0.00000000119'
Set-TextValue 'E47' '  +1.64%  '
Set-TextValue 'B48' 'EnergySwap'
Set-TextValue 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '9.228'
Set-TextValue 'E48' '  +0.93%  '
Set-TextValue 'D49' '0.4061'
Set-TextValue 'E49' '  +0.74%  '
Set-TextValue 'B50' 'RenderToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D50' '1.673'
Set-TextValue 'E50' '  +5.88%  '
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.05795'
Set-TextValue 'E51' '  +0.91%  '
